$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix J3 instruction text (was a lone space placeholder) ---
$ws.Range("J3").Value = "Bitte klicke an."

# --- Duplicate "SubjectArea" (AR) into the new "DescriptionOfPartialQualification" (AT)
#     column for every data row (2-56), matching the source column value per row ---
for ($r = 2; $r -le 56; $r++) {
    $arText = $ws.Cells.Item($r, 44).Value2
    $ws.Cells.Item($r, 46).Value = $arText
}

# --- Add new sample/course/survey metadata on row 2 ---
$ws.Cells.Item(2, 58).Value = 20                                     # BF2 Duration
$ws.Cells.Item(2, 61).Value = "Bertelsmann Stiftung"                 # BI2 Publisher
$ws.Cells.Item(2, 62).Value = "Selbsttest – Digitale Kompetenzen"    # BJ2 Title

# --- Update the selected cell to match the new focus point ---
$ws.Range("BI2").Select()
